# NYPD 107th Precinct CompStat weekly report refresh:
#  - New Police Commissioner name
#  - Volume/week number and report date range bumped to the next week
#  - Updated crime-complaint figures for the week/28-day/YTD/2yr/13yr/30yr columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header / masthead text updates
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# ---------------------------------------------------------------------------
# Helpers for the data table (rows 14-30).
#
# Two "placeholder" text styles are used throughout the table instead of a
# literal number when a figure is unavailable: text "0" (count columns) and
# text "***.*" (percent-change columns). When a cell flips between a real
# number and one of these text placeholders, the underlying style index also
# has to flip (s=14 for text, s=15/16 for numbers). Range.Copy() from a
# same-shaped template cell brings the right style+string along in one shot;
# we then stamp the real target value over the top when it is numeric.
# ---------------------------------------------------------------------------

function Set-NumCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-TextZero($addr) {
    # Stable template: C14 is the text "0" placeholder and is never touched.
    $ws.Range("C14").Copy($ws.Range($addr))
}

function Set-TextStar($addr) {
    # Stable template: E14 is the text "***.*" placeholder and is never touched.
    $ws.Range("E14").Copy($ws.Range($addr))
}

function Set-NumFromText($addr, $templateAddr, $val) {
    # Pulls a numeric style from a same-column template cell, then writes
    # the real numeric value.
    $ws.Range($templateAddr).Copy($ws.Range($addr))
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-NumCell "M15" 50

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-NumFromText "C16" "C17" 2
Set-NumCell "D16" 5
Set-NumCell "E16" -60
Set-NumCell "F16" 6
Set-NumCell "G16" 20
Set-NumCell "H16" -70
Set-NumCell "I16" 66
Set-NumCell "J16" 94
Set-NumCell "K16" -29.787234042553
Set-NumCell "L16" 26.923076923076
Set-NumCell "M16" -49.618320610687
Set-NumCell "N16" -86.221294363256

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-NumCell "C17" 5
Set-NumCell "D17" 4
Set-NumCell "E17" 25
Set-NumCell "F17" 15
Set-NumCell "G17" 20
Set-NumCell "H17" -25
Set-NumCell "I17" 110
Set-NumCell "J17" 96
Set-NumCell "K17" 14.583333333333
Set-NumCell "L17" 57.142857142857
Set-NumCell "M17" 129.166666666667
Set-NumCell "N17" -8.333333333333

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-NumCell "C18" 1
Set-NumCell "E18" -83.333333333333
Set-NumCell "G18" 27
Set-NumCell "H18" -51.851851851851
Set-NumCell "I18" 134
Set-NumCell "J18" 109
Set-NumCell "K18" 22.935779816513
Set-NumCell "L18" 42.553191489361
Set-NumCell "M18" 1.515151515151
Set-NumCell "N18" -82.908163265306

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-NumCell "C19" 13
Set-NumCell "D19" 10
Set-NumCell "E19" 30
Set-NumCell "F19" 43
Set-NumCell "G19" 63
Set-NumCell "H19" -31.746031746031
Set-NumCell "I19" 304
Set-NumCell "J19" 344
Set-NumCell "K19" -11.627906976744
Set-NumCell "L19" 92.405063291139
Set-NumCell "M19" 23.076923076923
Set-NumCell "N19" -5.882352941176

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-NumCell "C20" 12
Set-NumCell "D20" 5
Set-NumCell "E20" 140
Set-NumCell "G20" 19
Set-NumCell "H20" 152.631578947368
Set-NumCell "I20" 183
Set-NumCell "J20" 129
Set-NumCell "K20" 41.860465116279
Set-NumCell "L20" 221.052631578947
Set-NumCell "M20" 72.641509433962
Set-NumCell "N20" -92.868277474668

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold)
# ---------------------------------------------------------------------------
Set-NumCell "C21" 33
Set-NumCell "D21" 30
Set-NumCell "E21" 10
Set-NumCell "F21" 126
Set-NumCell "G21" 151
Set-NumCell "H21" -16.556291390728
Set-NumCell "I21" 806
Set-NumCell "J21" 785
Set-NumCell "K21" 2.675159235668
Set-NumCell "L21" 83.599088838268
Set-NumCell "M21" 19.762258543833
Set-NumCell "N21" -81.233993015133

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-NumCell "C22" 3
Set-TextZero "D22"
Set-TextStar "E22"
Set-NumCell "F22" 7
Set-NumCell "H22" 600
Set-NumCell "I22" 24
Set-NumCell "K22" 60
Set-NumCell "L22" 14.285714285714
Set-NumCell "M22" 71.428571428571

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
Set-NumCell "C23" 1
Set-NumCell "E23" 0
Set-NumCell "I23" 41
Set-NumCell "J23" 30
Set-NumCell "K23" 36.666666666666
Set-NumCell "L23" 127.777777777778
Set-NumCell "M23" 95.238095238095

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
Set-NumCell "C24" 13
Set-NumCell "D24" 38
Set-NumCell "E24" -65.789473684210
Set-NumCell "F24" 92
Set-NumCell "G24" 143
Set-NumCell "H24" -35.664335664335
Set-NumCell "I24" 748
Set-NumCell "J24" 712
Set-NumCell "K24" 5.056179775280
Set-NumCell "L24" 88.888888888888
Set-NumCell "M24" 49.003984063745

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
Set-NumCell "C25" 5
Set-NumCell "D25" 10
Set-NumCell "E25" -50
Set-NumCell "G25" 45
Set-NumCell "H25" -13.333333333333
Set-NumCell "I25" 244
Set-NumCell "J25" 242
Set-NumCell "K25" 0.826446280991
Set-NumCell "L25" 45.238095238095
Set-NumCell "M25" 12.962962962963

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-NumCell "F26" 1
Set-NumCell "H26" -75

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-NumFromText "C27" "C18" 2
Set-TextZero "D27"
Set-TextStar "E27"
Set-NumCell "F27" 3
Set-NumCell "G27" 3
Set-NumCell "H27" 0
Set-NumCell "I27" 24
Set-NumCell "K27" -29.411764705882
Set-NumCell "L27" -7.692307692307

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-NumFromText "D28" "D17" 1
Set-NumFromText "E28" "E17" -100
Set-NumFromText "G28" "G17" 1
Set-NumFromText "H28" "H17" -100
Set-NumCell "J28" 7
Set-NumCell "K28" -85.714285714285

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-NumFromText "D29" "D17" 1
Set-NumFromText "E29" "E17" -100
Set-NumFromText "G29" "G17" 1
Set-NumFromText "H29" "H17" -100
Set-NumCell "J29" 5
Set-NumCell "K29" -80

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
Set-NumFromText "C30" "C19" 1
Set-NumCell "F30" 3
Set-NumCell "I30" 16
Set-NumCell "K30" 700
Set-NumCell "L30" 433.333333333333
